# Ajuste do processo para executar apenas em uma unica chamada2
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap header labels for F1 / G1 ---
$ws.Range("F1").Value = "StatusEnvio"
$ws.Range("G1").Value = "Conexão"

# --- Row 2: replace with new PDV entry ---
$ws.Range("A2").Value = 5265
$ws.Range("B2").Value = "L5265 - SWIFT FLAMBOYANT (GO)"
$ws.Range("C2").Value = "PDV 01"
$ws.Range("D2").Value = "10.240.147.83"
$ws.Range("E2").Value = "Ping OK"
$ws.Range("F2").Value = "Enviado"
$ws.Range("G2").Value = "Conectado"

# --- Row 3: new PDV entry ---
$ws.Range("A3").Value = 5244
$ws.Range("B3").Value = "L5244 - SWIFT TAGUATINGA SUL (DF)"
$ws.Range("C3").Value = "PDV 03"
$ws.Range("D3").Value = "10.240.134.31"
$ws.Range("E3").Value = "Ping OK"
$ws.Range("F3").Value = "Enviado"
$ws.Range("G3").Value = "Conectado"

# --- Row 4: new PDV entry (E4 intentionally left blank) ---
$ws.Range("A4").Value = 5259
$ws.Range("B4").Value = "L5259 - SWIFT ALL WAYS MALL (RJ)"
$ws.Range("C4").Value = "PDV 02"
$ws.Range("D4").Value = "10.240.132.156"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Value = "Enviado"
$ws.Range("G4").Value = "Conectado"
